$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the split "PROBLEM DEFINITION: " runs (P / ROBLEM DEFINITION / ": ")
#    into a single run. The trailing space in the source is a non-breaking
#    space (U+00A0), preserve it exactly.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("PROBLEM DEFINITION: ", $false, $false, $false, $false, $false, $true, 1, $false, "PROBLEM DEFINITION: ", 2)

# ---------------------------------------------------------------------
# 2) Merge the "<label>:" + " " run pairs for the three numbered
#    objectives into single runs ending with a regular trailing space.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("1. Optimize User Experience (UX): ", $false, $false, $false, $false, $false, $true, 1, $false, "1. Optimize User Experience (UX): ", 2)

$rng = $d.Content
$rng.Find.Execute("2. Increase Conversion Rates: ", $false, $false, $false, $false, $false, $true, 1, $false, "2. Increase Conversion Rates: ", 2)

$rng = $d.Content
$rng.Find.Execute("3. Enhance Content Strategy: ", $false, $false, $false, $false, $false, $true, 1, $false, "3. Enhance Content Strategy: ", 2)

# ---------------------------------------------------------------------
# 3) Remove the "Design thinking can be a valuable approach..." intro
#    paragraph (and its following blank paragraph) right after the
#    "DESIGN THINKING:" heading, and remove the blank paragraph + the
#    "Remember, design thinking is an iterative process..." outro
#    paragraph right before "CONCLUSION:".
#
#    Paragraphs are located by their exact text and deleted from the
#    highest index down to the lowest so earlier indices stay valid
#    while later ones are removed.
# ---------------------------------------------------------------------
function Get-ParaIndexByText($doc, [string]$needle) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        if ($paras.Item($i).Range.Text.TrimEnd("`r") -eq $needle) {
            return $i
        }
    }
    return -1
}

# --- outro: blank paragraph + "Remember, design thinking..." ---
$idxRemember = Get-ParaIndexByText $d "Remember, design thinking is an iterative process, so you may need to revisit and refine your website traffic analysis dashboard as user needs and goals evolve."
$idxBlankBeforeRemember = $idxRemember - 1

$paras = $d.Paragraphs
$paras.Item($idxRemember).Range.Delete()
$paras = $d.Paragraphs
$paras.Item($idxBlankBeforeRemember).Range.Delete()

# --- intro: "Design thinking can be a valuable approach..." + blank paragraph ---
$idxIntro = Get-ParaIndexByText $d "Design thinking can be a valuable approach for website traffic analysis to ensure you're addressing user needs effectively. Here's a simplified process:"
$idxBlankAfterIntro = $idxIntro + 1

$paras = $d.Paragraphs
$paras.Item($idxBlankAfterIntro).Range.Delete()
$paras = $d.Paragraphs
$paras.Item($idxIntro).Range.Delete()

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
